$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1288.7142
$ws.Range("I9").Value = 1319.4546
$ws.Range("K9").Value = 1319.4546
$ws.Range("M9").Value = -1150.4546
$ws.Range("H15").Value = 2395.7046
$ws.Range("I15").Value = 2395.7046
$ws.Range("K15").Value = 7187.1138
$ws.Range("M15").Value = -7018.1138
$ws.Range("H40").Value = 1090
$ws.Range("J40").Value = 978.5
$ws.Range("L40").Value = 978.5
$ws.Range("N40").Value = -1328.5
$ws.Range("H58").Value = 12657.111
$ws.Range("J58").Value = 4112.5
$ws.Range("L58").Value = 12337.5
$ws.Range("N58").Value = -12637.5
$ws.Range("H70").Value = 6353.6665
$ws.Range("I70").Value = 2809.0908
$ws.Range("J70").Value = 16101.25
$ws.Range("K70").Value = 8427.2724
$ws.Range("L70").Value = 48303.75
$ws.Range("M70").Value = -8157.2724
$ws.Range("N70").Value = -48843.75
$ws.Range("H73").Value = 6353.6665
$ws.Range("I73").Value = 2809.0908
$ws.Range("J73").Value = 16101.25
$ws.Range("K73").Value = 8427.2724
$ws.Range("L73").Value = 48303.75
$ws.Range("M73").Value = -7491.2724
$ws.Range("N73").Value = -50175.75
$ws.Range("H112").Value = 5375.304
$ws.Range("J112").Value = 3295.5334
$ws.Range("L112").Value = 9886.600199999999
$ws.Range("N112").Value = -12102.6002
$ws.Range("H113").Value = 14178
$ws.Range("I113").Value = 18850.4
$ws.Range("J113").Value = 2497
$ws.Range("K113").Value = 18850.4
$ws.Range("L113").Value = 2497
$ws.Range("M113").Value = -15596.4
$ws.Range("N113").Value = -9005
$ws.Range("H137").Value = 8259.487999999999
$ws.Range("I137").Value = 2952.1667
$ws.Range("K137").Value = 8856.500100000001
$ws.Range("M137").Value = -6306.500100000001
$ws.Range("H138").Value = 3698.678
$ws.Range("I138").Value = 4596.875
$ws.Range("J138").Value = 3557.7844
$ws.Range("K138").Value = 13790.625
$ws.Range("L138").Value = 10673.3532
$ws.Range("M138").Value = -8650.625
$ws.Range("N138").Value = -20953.3532

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 5817.1665
$ws.Range("I36").Value = 1999.5
$ws.Range("J36").Value = 7726
$ws.Range("K36").Value = 1999.5
$ws.Range("L36").Value = 7726
$ws.Range("M36").Value = -1653.5
$ws.Range("N36").Value = -8418
$ws.Range("H45").Value = 2349.6924
$ws.Range("I45").Value = 2249.6365
$ws.Range("J45").Value = 2900
$ws.Range("K45").Value = 2249.6365
$ws.Range("L45").Value = 2900
$ws.Range("M45").Value = -1872.6365
$ws.Range("N45").Value = -3654
$ws.Range("H97").Value = 4282.0454
$ws.Range("I97").Value = 2160.3
$ws.Range("K97").Value = 2160.3
$ws.Range("M97").Value = -1664.3
$ws.Range("H110").Value = 4513.269
$ws.Range("I110").Value = 2936.8572
$ws.Range("K110").Value = 2936.8572
$ws.Range("M110").Value = -891.8571999999999
$ws.Range("H122").Value = 668366.5
$ws.Range("I122").Value = 770676.9399999999
$ws.Range("J122").Value = 3348.75
$ws.Range("K122").Value = 2312030.82
$ws.Range("L122").Value = 10046.25
$ws.Range("M122").Value = -2309580.82
$ws.Range("N122").Value = -14946.25
$ws.Range("H132").Value = 2457126.5
$ws.Range("I132").Value = 4690.814
$ws.Range("K132").Value = 14072.442
$ws.Range("M132").Value = -11542.442

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 28773.104
$ws.Range("J20").Value = 38538.473
$ws.Range("L20").Value = 38538.473
$ws.Range("N20").Value = -39032.473
$ws.Range("H22").Value = 3999.75
$ws.Range("I22").Value = 3999.75
$ws.Range("K22").Value = 3999.75
$ws.Range("M22").Value = -3826.75
$ws.Range("H99").Value = 11698.75
$ws.Range("I99").Value = 15920.5
$ws.Range("J99").Value = 1848
$ws.Range("K99").Value = 15920.5
$ws.Range("L99").Value = 1848
$ws.Range("M99").Value = -14422.5
$ws.Range("N99").Value = -4844
$ws.Range("H107").Value = 600.3
$ws.Range("I107").Value = 550.625
$ws.Range("K107").Value = 550.625
$ws.Range("M107").Value = 1369.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1277
$ws.Range("I22").Value = 694.5
$ws.Range("J22").Value = 1665.3334
$ws.Range("K22").Value = 694.5
$ws.Range("L22").Value = 1665.3334
$ws.Range("M22").Value = -344.5
$ws.Range("N22").Value = -2365.3334
$ws.Range("H31").Value = 75168.586
$ws.Range("I31").Value = 157718.69
$ws.Range("J31").Value = 24066.143
$ws.Range("K31").Value = 157718.69
$ws.Range("L31").Value = 24066.143
$ws.Range("M31").Value = -157423.69
$ws.Range("N31").Value = -24656.143
$ws.Range("H34").Value = 75168.586
$ws.Range("I34").Value = 157718.69
$ws.Range("J34").Value = 24066.143
$ws.Range("K34").Value = 157718.69
$ws.Range("L34").Value = 24066.143
$ws.Range("M34").Value = -157516.69
$ws.Range("N34").Value = -24470.143
$ws.Range("H58").Value = 22328.2
$ws.Range("I58").Value = 10073.375
$ws.Range("J58").Value = 30498.084
$ws.Range("K58").Value = 10073.375
$ws.Range("L58").Value = 30498.084
$ws.Range("M58").Value = -9870.375
$ws.Range("N58").Value = -30904.084
$ws.Range("H99").Value = 7706.567
$ws.Range("I99").Value = 5145.8125
$ws.Range("K99").Value = 5145.8125
$ws.Range("M99").Value = -3647.8125
$ws.Range("H126").Value = 7706.567
$ws.Range("I126").Value = 5145.8125
$ws.Range("K126").Value = 15437.4375
$ws.Range("M126").Value = -12967.4375
$ws.Range("H132").Value = 41184270
$ws.Range("I132").Value = 3059
$ws.Range("J132").Value = 102956090
$ws.Range("K132").Value = 9177
$ws.Range("L132").Value = 308868270
$ws.Range("M132").Value = -6647
$ws.Range("N132").Value = -308873330
$ws.Range("H136").Value = 22328.2
$ws.Range("I136").Value = 10073.375
$ws.Range("J136").Value = 30498.084
$ws.Range("K136").Value = 30220.125
$ws.Range("L136").Value = 91494.25199999999
$ws.Range("M136").Value = -27670.125
$ws.Range("N136").Value = -96594.25199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 25.6
$ws.Range("J2").Value = 9.800000000000001
$ws.Range("L2").Value = 58.8
$ws.Range("N2").Value = -284.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 20987.223
$ws.Range("I70").Value = 12777.2
$ws.Range("K70").Value = 12777.2
$ws.Range("M70").Value = -12507.2
$ws.Range("H73").Value = 20987.223
$ws.Range("I73").Value = 12777.2
$ws.Range("K73").Value = 12777.2
$ws.Range("M73").Value = -11841.2
$ws.Range("H97").Value = 4642.2144
$ws.Range("I97").Value = 1454.5714
$ws.Range("J97").Value = 7829.857
$ws.Range("K97").Value = 1454.5714
$ws.Range("L97").Value = 7829.857
$ws.Range("M97").Value = -958.5714
$ws.Range("N97").Value = -8821.857
$ws.Range("H102").Value = 7690.0557
$ws.Range("I102").Value = 5432.6875
$ws.Range("K102").Value = 5432.6875
$ws.Range("M102").Value = -3810.6875
$ws.Range("H122").Value = 7254.1763
$ws.Range("I122").Value = 4016.5
$ws.Range("K122").Value = 12049.5
$ws.Range("M122").Value = -9599.5
$ws.Range("H132").Value = 392130.28
$ws.Range("I132").Value = 5309.3105
$ws.Range("J132").Value = 1139984.1
$ws.Range("K132").Value = 15927.9315
$ws.Range("L132").Value = 3419952.3
$ws.Range("M132").Value = -13397.9315
$ws.Range("N132").Value = -3425012.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 13634.4
$ws.Range("I7").Value = 14968.9
$ws.Range("K7").Value = 14968.9
$ws.Range("M7").Value = -14856.9
$ws.Range("H9").Value = 1868.6
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 1868.6
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 1868.6
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -2316.6
$ws.Range("H40").Value = 8863.143
$ws.Range("I40").Value = 6044
$ws.Range("J40").Value = 12622
$ws.Range("K40").Value = 6044
$ws.Range("L40").Value = 12622
$ws.Range("M40").Value = -5908
$ws.Range("N40").Value = -12894
$ws.Range("H46").Value = 2173.8333
$ws.Range("I46").Value = 986.75
$ws.Range("K46").Value = 986.75
$ws.Range("M46").Value = -798.75
$ws.Range("H126").Value = 13634.4
$ws.Range("I126").Value = 14968.9
$ws.Range("K126").Value = 44906.7
$ws.Range("M126").Value = -42436.7
$ws.Range("H132").Value = 1948572.2
$ws.Range("I132").Value = 7875.2144
$ws.Range("J132").Value = 3183561.2
$ws.Range("K132").Value = 23625.6432
$ws.Range("L132").Value = 9550683.600000001
$ws.Range("M132").Value = -21095.6432
$ws.Range("N132").Value = -9555743.600000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H81").Value = 6680.8184
$ws.Range("J81").Value = 6783.3335
$ws.Range("L81").Value = 13566.667
$ws.Range("N81").Value = -15688.667
$ws.Range("H84").Value = 6680.8184
$ws.Range("J84").Value = 6783.3335
$ws.Range("L84").Value = 67833.33499999999
$ws.Range("N84").Value = -78441.33499999999
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 246832.58
$ws.Range("I136").Value = 2781.9333
$ws.Range("K136").Value = 8345.7999
$ws.Range("M136").Value = -5795.7999
